$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("surveys")

# Insert a new column before column D (fieldPeriod.start), shifting everything right.
$ws.Columns("D").EntireColumn.Insert()

# New column D header + values ("wave")
$ws.Range("D1").Value = "wave"
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 1

# New column D width (closest representable value under this engine's
# 1/6-character rounding of ColumnWidth; target stored width is 15.7109375)
$ws.Columns("D").ColumnWidth = 14.833333333333334

# Column B width changes (title.de); target stored width is 33.5703125
$ws.Columns("B").ColumnWidth = 32.666666666666664

# Apply wrap text formatting to B1:C3 (title.de / title.en columns)
$ws.Range("B1:C3").WrapText = $true

# Update the view/selection
$ws.Range("D4").Select()

Write-Host "done"
